# Applies the CodeSystem-claim-supporting-info-category.xlsx metadata update:
#  - Version bumped 5.0.0 -> 6.0.0
#  - Date updated
#  - Publisher value filled in ("Alvearie Team")
#  - Duplicate "Contact" row replaced by a "Jurisdiction" row, and the
#    leftover duplicate "Contact" row removed (sheet shrinks from 22 to 21 rows)
#  - "Case Sensitive" value filled in ("true")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$ws.Range("B3").Value = "6.0.0"

# Update Date value (row 8)
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Fill in Publisher value (row 9)
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 used to be "Contact" / "No display for ContactDetail" and row 11
# duplicated it. Turn row 10 into the new "Jurisdiction" row and delete the
# now-redundant row 11, shifting everything below up by one.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
$ws.Rows.Item(11).Delete()

# Fill in Case Sensitive value (now row 14 after the row deletion).
# A bare "true" would be auto-coerced to a boolean cell by the engine, so
# instead enter it as a formula that evaluates to the text string "true"
# and then convert the cell to a plain value in place; this yields a
# literal text cell without disturbing the existing cell formatting/style.
$ws.Range("B14").Formula = "=""true"""
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0
